$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("key_outputs")

# Remove the columns G and H values for the data rows (2-8); clear contents
# so the cells no longer exist in the saved sheet.
$ws.Range("G2:H8").ClearContents()

# Update the active selection on the sheet to F2 (was F6).
$ws.Range("F2").Select()
